$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value (Excel serial date) for every
# data row (2 through 344). The value was bumped from 45179 (2023-09-10)
# to 45180 (2023-09-11) for all of these rows.
$ws.Range("C2:C344").Value = 45180
